# Spring Boot deck: insert "Agenda" / "Spring Boot" / "We want more!" slides
# right after "About us", and a "LAB 3" slide between "LAB 2" and "Q&A".
#
# Slides are created in the same order PowerPoint assigned their internal
# SlideIDs (Agenda, Spring Boot, LAB 3, We want more!) so each lands on its
# final index immediately - no extra MoveTo shuffling required:
#   1 About/title slides stay put
#   3 Agenda            (new)
#   4 Spring Boot        (new)
#   5 We want more!      (new, inserted after Spring Boot once LAB 3 exists)
#   6-10 existing Acme architecture .. LAB 2 slides shift down automatically
#   11 LAB 3: centralize messagehandling (optional)  (new, right before Q&A)
#   12 Q&A

$p = $ppt.ActivePresentation

$master = $p.Designs.Item(1).SlideMaster
$layoutTitleContent = $master.CustomLayouts.Item(2)   # title + content placeholder layout

function Set-NoBullet($textRange, [int]$paraIndex) {
    $para = $textRange.Paragraphs($paraIndex, 1)
    $para.ParagraphFormat.Bullet.Type = 0
}

# ---------------------------------------------------------------------
# "Agenda" (lands at slide 3)
# ---------------------------------------------------------------------
$sAgenda = $p.Slides.AddSlide(3, $layoutTitleContent)
$sAgenda.Shapes.Item(1).Name = "Titel 1"
$sAgenda.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"
$sAgenda.Shapes.Item(1).TextFrame.TextRange.Text = "Agenda"

$bodyAgenda = $sAgenda.Shapes.Item(2).TextFrame.TextRange
$bodyAgenda.Text = "Spring Boot`rWe want more!`r`rLab1`rLab2`r(Lab3)`r`rWrap up"

# ---------------------------------------------------------------------
# "Spring Boot" mission-statement bullets (lands at slide 4)
# ---------------------------------------------------------------------
$sSpringBoot = $p.Slides.AddSlide(4, $layoutTitleContent)
$sSpringBoot.Shapes.Item(1).Name = "Titel 1"
$sSpringBoot.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"
$sSpringBoot.Shapes.Item(1).TextFrame.TextRange.Text = "Spring Boot"

$bodySpringBoot = $sSpringBoot.Shapes.Item(2).TextFrame.TextRange
$bodySpringBoot.Text = "Provide a radically faster and widely accessible getting started experience for all Spring development.`rBe opinionated out of the box, but get out of the way quickly as requirements start to diverge from the defaults.`rProvide a range of non-functional features that are common to large classes of projects (e.g. embedded servers, security, metrics, health checks, externalized configuration).`rAbsolutely no code generation and no requirement for XML configuration.`r"

# ---------------------------------------------------------------------
# "LAB 3: centralize messagehandling (optional)" (lands right before Q&A,
# currently the last slide - create it now so SlideIDs stay in sequence).
# ---------------------------------------------------------------------
$qaIndex = $p.Slides.Count
$sLab3 = $p.Slides.AddSlide($qaIndex, $layoutTitleContent)
$sLab3.Shapes.Item(1).Name = "Titel 1"
$sLab3.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"
$sLab3.Shapes.Item(1).TextFrame.TextRange.Text = "LAB 3: centralize messagehandling (optional)"

$bodyLab3 = $sLab3.Shapes.Item(2).TextFrame.TextRange
$bodyLab3.Text = "Step 1: Remove any application-specific implementation of JMS messagelisteners`r`rStep 2: Add the @MessageHandlerBean annotion to the messagehandler"

# ---------------------------------------------------------------------
# "We want more!" (lands at slide 5, right after "Spring Boot")
# ---------------------------------------------------------------------
$sWeWant = $p.Slides.AddSlide(5, $layoutTitleContent)
$sWeWant.Shapes.Item(1).Name = "Titel 1"
$sWeWant.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"
$sWeWant.Shapes.Item(1).TextFrame.TextRange.Text = "We want more!"

$bodyWeWant = $sWeWant.Shapes.Item(2).TextFrame.TextRange
$bodyWeWant.Text = "Imagine your company:`rhas a different opinion`rbuilds software on a uniform stack`rhas a lot of applications to maintain`rhas custom requirements across all applications"
Set-NoBullet $bodyWeWant 1

# ---------------------------------------------------------------------
# Refresh the cached "today" date shown on the Notes / Handout masters.
# ---------------------------------------------------------------------
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "05/15/2017"
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = "05/15/2017"

Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    Write-Output "$i -> id=$($sl.SlideID) title=$($sl.Shapes.Item(1).TextFrame.TextRange.Text)"
}
